# Apply the edits described in the diff to the Tab17 workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab17")

# --- Fix "Etats" -> "États" typos and add missing "RDM, pays en developpement sans littoral" label ---
$ws.Range("B93").Value = "Afrique, petits États insulaires en développement"
$ws.Range("B94").Value = "RDM, petits États insulaires en développement"
$ws.Range("B96").Value = "RDM, pays en développement sans littoral"
$ws.Range("B97").Value = "Afrique, États fragiles"
$ws.Range("B98").Value = "RDM, États fragiles"

# --- Corrected data values in row 92 (RDM, pays les moins avances) ---
$ws.Range("C92").Value = 76.768915669000805
$ws.Range("D92").Value = 12.751171563349899
$ws.Range("G92").Value = 22.719498363479499
$ws.Range("H92").Value = 18.5495027374903
$ws.Range("I92").Value = 4.1699326226786901

# --- Update the responsibility disclaimer text ---
$ws.Range("A104").Value = "Responsabilité : Ce tableau ainsi que toutes les données qu'il peut comprendre, sont sans préjudice du statut de tout territoire, de la souveraineté s'exerçant sur ce dernier, du tracé des frontières et limites internationales, et du nom de tout territoire, ville ou région."

# --- Window size in the workbook view ---
$excel.ActiveWindow.Width = 19200
$excel.ActiveWindow.Height = 5590
